# "Egg stir fried rice." - add a missing "ground cloves" ingredient row
# to the Ingredients sheet (inserted at row 83, pushing the rest down).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 83 (everything from the old row 83 onward
# shifts down by one, exactly like pressing Ctrl+"+" on a selected row).
$ws.Rows("83:83").Insert()

# Populate the new row with the new ingredient.
$ws.Cells.Item(83, 1).Value = "ground cloves"
$ws.Cells.Item(83, 2).Value = "Check"
$ws.Cells.Item(83, 3).Value = 0
$ws.Cells.Item(83, 4).Value = 0
$ws.Cells.Item(83, 5).Value = 0
$ws.Cells.Item(83, 6).Value = 0

# The worksheet's AutoFilter range grew by one row.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Ingredients!_FilterDatabase") {
        $n.RefersTo = "=Ingredients!`$A`$2:`$G`$175"
    }
}

# Reflect where the user scrolled to / what they selected while typing
# the new row in: frozen header row, scrolled down near the new data,
# with B83:F83 selected.
$excel.ActiveWindow.FreezePanes = $false
$ws.Range("A74").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("B83:F83").Select()
